$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '85.773.12'
$ws.Range('E2').Value = '  +4.83%  '

# Row 3
Set-TextValue 'D3' '3.250.82'
$ws.Range('E3').Value = '  +2.68%  '

# Row 4
$ws.Range('E4').Value = '  +0.58%  '

# Row 5
Set-TextValue 'D5' '208.29'
$ws.Range('E5').Value = '  -4.38%  '

# Row 6
Set-TextValue 'D6' '619.92'
$ws.Range('E6').Value = '  +0.27%  '

# Row 7
Set-TextValue 'D7' '0.357'
$ws.Range('E7').Value = '  +24.21%  '

# Row 8
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue 'D8' '0.648'
$ws.Range('E8').Value = '  +11.89%  '

# Row 9
$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue 'D9' '1.00'
$ws.Range('E9').Value = '  +0.26%  '

# Row 10
Set-TextValue 'D10' '3.262.28'
$ws.Range('E10').Value = '  +3.34%  '

# Row 11
Set-TextValue 'D11' '0.571'
$ws.Range('E11').Value = '  -2.45%  '

# Row 12
$ws.Range('E12').Value = '  +6.77%  '

# Row 13
Set-TextValue 'D13' '0.0000252'
$ws.Range('E13').Value = '  -0.51%  '

# Row 14
Set-TextValue 'D14' '3.864.55'
$ws.Range('E14').Value = '  +3.27%  '

# Row 15
Set-TextValue 'D15' '33.60'
$ws.Range('E15').Value = '  +5.27%  '

# Row 16
Set-TextValue 'D16' '5.25'
$ws.Range('E16').Value = '  -0.82%  '

# Row 17
Set-TextValue 'D17' '85.633.69'
$ws.Range('E17').Value = '  +4.79%  '

# Row 18
Set-TextValue 'D18' '3.257.13'
$ws.Range('E18').Value = '  +2.82%  '

# Row 19
Set-TextValue 'D19' '13.94'
$ws.Range('E19').Value = '  +0.28%  '

# Row 20
$ws.Range('E20').Value = '  -7.95%  '

# Row 21
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D21' '425.59'
$ws.Range('E21').Value = '  -2.31%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D22' '8.88'
$ws.Range('E22').Value = '  +0.34%  '

# Row 23
Set-TextValue 'D23' '5.26'
$ws.Range('E23').Value = '  +2.84%  '

# Row 24
Set-TextValue 'D24' '7.20'
$ws.Range('E24').Value = '  -0.85%  '

# Row 25
Set-TextValue 'D25' '12.35'
$ws.Range('E25').Value = '  +4.59%  '

# Row 26
Set-TextValue 'D26' '5.05'
$ws.Range('E26').Value = '  -3.09%  '

# Row 27
Set-TextValue 'D27' '3.425.09'
$ws.Range('E27').Value = '  +3.36%  '

# Row 28
Set-TextValue 'D28' '75.20'
$ws.Range('E28').Value = '  -1.44%  '

# Row 29
$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D29' '1.00'
$ws.Range('E29').Value = '  +0.16%  '

# Row 30
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D30' '0.0000127'
$ws.Range('E30').Value = '  +5.58%  '

# Row 31
$ws.Range('E31').Value = '  +17.67%  '

# Row 32
$ws.Range('E32').Value = '  +0.16%  '

# Row 33
$ws.Range('E33').Value = '  -2.40%  '

# Row 34
Set-TextValue 'D34' '541.75'
$ws.Range('E34').Value = '  -4.16%  '

# Row 35
$ws.Range('E35').Value = '  -4.58%  '

# Row 36
Set-TextValue 'D36' '1.93'
$ws.Range('E36').Value = '  -2.34%  '

# Row 37
Set-TextValue 'D37' '6.82'
$ws.Range('E37').Value = '  +10.69%  '

# Row 38
Set-TextValue 'D38' '0.136'
$ws.Range('E38').Value = '  -9.50%  '

# Row 39
Set-TextValue 'D39' '22.24'
$ws.Range('E39').Value = '  -1.21%  '

# Row 40
Set-TextValue 'D40' '1.01'
$ws.Range('E40').Value = '  +0.82%  '

# Row 41
Set-TextValue 'D41' '21.57'
$ws.Range('E41').Value = '  +3.56%  '

# Row 42
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue 'D42' '0.388'
$ws.Range('E42').Value = '  -3.61%  '

# Row 43
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D43' '1.97'
$ws.Range('E43').Value = '  -1.35%  '

# Row 44
Set-TextValue 'D44' '158.16'
$ws.Range('E44').Value = '  -0.38%  '

# Row 45
$ws.Range('E45').Value = '  -0.06%  '

# Row 46
$ws.Range('E46').Value = '  -3.33%  '

# Row 47
Set-TextValue 'D47' '177.01'
$ws.Range('E47').Value = '  -4.70%  '

# Row 48
Set-TextValue 'D48' '44.00'
$ws.Range('E48').Value = '  -0.97%  '

# Row 49
Set-TextValue 'D49' '1.28'
$ws.Range('E49').Value = '  -1.89%  '

# Row 50
Set-TextValue 'D50' '4.21'
$ws.Range('E50').Value = '  +1.32%  '

# Row 51
Set-TextValue 'D51' '0.728'
$ws.Range('E51').Value = '  -4.25%  '
